$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was updated
# from 45172 (2023-09-03) to 45175 (2023-09-06) for every data row (2-89).
for ($row = 2; $row -le 89; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
